# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 6, pushing the existing rows 6-12
# down to rows 7-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6; Excel shifts rows 6:12 down to 7:13
# and carries formatting (e.g. the date style on column D) along with them.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with this week's data.
$ws.Cells.Item(6, 1).Value = 7
$ws.Cells.Item(6, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(6, 3).Value = "Ñuble"
$ws.Cells.Item(6, 4).Value = 44894
$ws.Cells.Item(6, 5).Value = 16
$ws.Cells.Item(6, 6).Value = 100114007
$ws.Cells.Item(6, 7).Value = "Jengibre"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 30
$ws.Cells.Item(6, 11).Value = 18000
$ws.Cells.Item(6, 12).Value = 18000
$ws.Cells.Item(6, 13).Value = 18000
$ws.Cells.Item(6, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(6, 15).Value = "Perú"
$ws.Cells.Item(6, 16).Value = 1385
$ws.Cells.Item(6, 17).Value = 13
$ws.Cells.Item(6, 18).Value = "Hortaliza"
